$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendor_payment")

# Update the Reference / Payment reference values in row 2 (H2 and I2)
# from "CR416232" to "CR414626"
$ws.Range("H2").Value = "CR414626"
$ws.Range("I2").Value = "CR414626"
